$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 2635.3777
$ws.Cells.Item(138, 9).Value = 1766.1428
$ws.Cells.Item(138, 10).Value = 5677.7
$ws.Cells.Item(138, 11).Value = 5298.428400000001
$ws.Cells.Item(138, 12).Value = 17033.1
$ws.Cells.Item(138, 13).Value = -158.4284000000007
$ws.Cells.Item(138, 14).Value = -27313.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 4549165
$ws.Cells.Item(74, 9).Value = 6452369.5
$ws.Cells.Item(74, 10).Value = 10753.77
$ws.Cells.Item(74, 11).Value = 6452369.5
$ws.Cells.Item(74, 12).Value = 10753.77
$ws.Cells.Item(74, 13).Value = -6451495.5
$ws.Cells.Item(74, 14).Value = -12501.77

$ws.Cells.Item(77, 8).Value = 4549165
$ws.Cells.Item(77, 9).Value = 6452369.5
$ws.Cells.Item(77, 10).Value = 10753.77
$ws.Cells.Item(77, 11).Value = 32261847.5
$ws.Cells.Item(77, 12).Value = 53768.85000000001
$ws.Cells.Item(77, 13).Value = -32257479.5
$ws.Cells.Item(77, 14).Value = -62504.85000000001

$ws.Cells.Item(132, 8).Value = 46243.914
$ws.Cells.Item(132, 9).Value = 2902.2666
$ws.Cells.Item(132, 10).Value = 127509.5
$ws.Cells.Item(132, 11).Value = 8706.799800000001
$ws.Cells.Item(132, 12).Value = 382528.5
$ws.Cells.Item(132, 13).Value = -6176.799800000001
$ws.Cells.Item(132, 14).Value = -387588.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1273.4286
$ws.Cells.Item(20, 9).Value = 1347.6364
$ws.Cells.Item(20, 10).Value = 1191.8
$ws.Cells.Item(20, 11).Value = 1347.6364
$ws.Cells.Item(20, 12).Value = 1191.8
$ws.Cells.Item(20, 13).Value = -1100.6364
$ws.Cells.Item(20, 14).Value = -1685.8

$ws.Cells.Item(82, 8).Value = 12675.923
$ws.Cells.Item(82, 9).Value = 5902.625
$ws.Cells.Item(82, 10).Value = 23513.2
$ws.Cells.Item(82, 11).Value = 5902.625
$ws.Cells.Item(82, 12).Value = 23513.2
$ws.Cells.Item(82, 13).Value = -5519.625
$ws.Cells.Item(82, 14).Value = -24279.2

$ws.Cells.Item(85, 8).Value = 12675.923
$ws.Cells.Item(85, 9).Value = 5902.625
$ws.Cells.Item(85, 10).Value = 23513.2
$ws.Cells.Item(85, 11).Value = 5902.625
$ws.Cells.Item(85, 12).Value = 23513.2
$ws.Cells.Item(85, 13).Value = -4576.625
$ws.Cells.Item(85, 14).Value = -26165.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 9).Value = 1169.15
$ws.Cells.Item(31, 10).Value = 2017.6666
$ws.Cells.Item(31, 11).Value = 1169.15
$ws.Cells.Item(31, 12).Value = 2017.6666
$ws.Cells.Item(31, 13).Value = -874.1500000000001
$ws.Cells.Item(31, 14).Value = -2607.6666

$ws.Cells.Item(34, 9).Value = 1169.15
$ws.Cells.Item(34, 10).Value = 2017.6666
$ws.Cells.Item(34, 11).Value = 1169.15
$ws.Cells.Item(34, 12).Value = 2017.6666
$ws.Cells.Item(34, 13).Value = -967.1500000000001
$ws.Cells.Item(34, 14).Value = -2421.6666

$ws.Cells.Item(68, 8).Value = 20450
$ws.Cells.Item(68, 10).Value = 20450
$ws.Cells.Item(68, 12).Value = 20450
$ws.Cells.Item(68, 14).Value = -21948

$ws.Cells.Item(71, 8).Value = 20450
$ws.Cells.Item(71, 10).Value = 20450
$ws.Cells.Item(71, 12).Value = 61350
$ws.Cells.Item(71, 14).Value = -68838

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(23, 8).Value = 132
$ws.Cells.Item(23, 9).Value = 43.333332
$ws.Cells.Item(23, 10).Value = 198.5
$ws.Cells.Item(23, 11).Value = 129.999996
$ws.Cells.Item(23, 12).Value = 595.5
$ws.Cells.Item(23, 13).Value = 105.000004
$ws.Cells.Item(23, 14).Value = -1065.5

$ws.Cells.Item(36, 8).Value = 400.66666
$ws.Cells.Item(36, 9).Value = 101
$ws.Cells.Item(36, 10).Value = 1000
$ws.Cells.Item(36, 11).Value = 303
$ws.Cells.Item(36, 12).Value = 3000
$ws.Cells.Item(36, 13).Value = -134
$ws.Cells.Item(36, 14).Value = -3338

$ws.Cells.Item(37, 8).Value = 40221.89
$ws.Cells.Item(37, 10).Value = 40221.89
$ws.Cells.Item(37, 12).Value = 120665.67
$ws.Cells.Item(37, 14).Value = -120889.67

$ws.Cells.Item(54, 8).Value = 3000
$ws.Cells.Item(54, 10).Value = 3000
$ws.Cells.Item(54, 12).Value = 9000
$ws.Cells.Item(54, 14).Value = -10118

$ws.Cells.Item(94, 8).Value = 3349.5
$ws.Cells.Item(94, 10).Value = 4024.6667
$ws.Cells.Item(94, 12).Value = 12074.0001
$ws.Cells.Item(94, 14).Value = -13426.0001

$ws.Cells.Item(95, 8).Value = 2500
$ws.Cells.Item(95, 10).Value = 2500
$ws.Cells.Item(95, 12).Value = 7500
$ws.Cells.Item(95, 14).Value = -11618

$ws.Cells.Item(99, 8).Value = 25
$ws.Cells.Item(99, 9).Value = 25
$ws.Cells.Item(99, 10).Value = 0
$ws.Cells.Item(99, 11).Value = 75
$ws.Cells.Item(99, 12).Value = 0
$ws.Cells.Item(99, 13).Value = 2171
$ws.Cells.Item(99, 14).ClearContents()

$ws.Cells.Item(100, 8).Value = 3112.8
$ws.Cells.Item(100, 10).Value = 4028
$ws.Cells.Item(100, 12).Value = 12084
$ws.Cells.Item(100, 14).Value = -13706

$ws.Cells.Item(101, 8).Value = 5716.8
$ws.Cells.Item(101, 10).Value = 6839.5
$ws.Cells.Item(101, 12).Value = 20518.5
$ws.Cells.Item(101, 14).Value = -25386.5

$ws.Cells.Item(105, 8).Value = 129716856
$ws.Cells.Item(105, 10).Value = 129716856
$ws.Cells.Item(105, 12).Value = 389150568
$ws.Cells.Item(105, 14).Value = -389155810

$ws.Cells.Item(108, 8).Value = 3030
$ws.Cells.Item(108, 10).Value = 3030
$ws.Cells.Item(108, 12).Value = 9090
$ws.Cells.Item(108, 14).Value = -14850

$ws.Cells.Item(110, 8).Value = 2572.111
$ws.Cells.Item(110, 9).Value = 2249.8572
$ws.Cells.Item(110, 10).Value = 3700
$ws.Cells.Item(110, 11).Value = 6749.571599999999
$ws.Cells.Item(110, 12).Value = 11100
$ws.Cells.Item(110, 13).Value = -2659.571599999999
$ws.Cells.Item(110, 14).Value = -19280

$ws.Cells.Item(131, 8).Value = 3587.0833
$ws.Cells.Item(131, 10).Value = 2669.4644
$ws.Cells.Item(131, 12).Value = 8008.3932
$ws.Cells.Item(131, 14).Value = -18088.3932

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 4543.6924
$ws.Cells.Item(70, 9).Value = 4466.8
$ws.Cells.Item(70, 10).Value = 4800
$ws.Cells.Item(70, 11).Value = 4466.8
$ws.Cells.Item(70, 12).Value = 4800
$ws.Cells.Item(70, 13).Value = -4196.8
$ws.Cells.Item(70, 14).Value = -5340

$ws.Cells.Item(73, 8).Value = 4543.6924
$ws.Cells.Item(73, 9).Value = 4466.8
$ws.Cells.Item(73, 10).Value = 4800
$ws.Cells.Item(73, 11).Value = 4466.8
$ws.Cells.Item(73, 12).Value = 4800
$ws.Cells.Item(73, 13).Value = -3530.8
$ws.Cells.Item(73, 14).Value = -6672

$ws.Cells.Item(132, 8).Value = 2493.3215
$ws.Cells.Item(132, 9).Value = 2433.8096
$ws.Cells.Item(132, 10).Value = 2671.8572
$ws.Cells.Item(132, 11).Value = 7301.4288
$ws.Cells.Item(132, 12).Value = 8015.571599999999
$ws.Cells.Item(132, 13).Value = -4771.4288
$ws.Cells.Item(132, 14).Value = -13075.5716

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 518.2857
$ws.Cells.Item(22, 9).Value = 345
$ws.Cells.Item(22, 10).Value = 587.6
$ws.Cells.Item(22, 11).Value = 345
$ws.Cells.Item(22, 12).Value = 587.6
$ws.Cells.Item(22, 13).Value = -50
$ws.Cells.Item(22, 14).Value = -1177.6

$ws.Cells.Item(27, 8).Value = 518.2857
$ws.Cells.Item(27, 9).Value = 345
$ws.Cells.Item(27, 10).Value = 587.6
$ws.Cells.Item(27, 11).Value = 345
$ws.Cells.Item(27, 12).Value = 587.6
$ws.Cells.Item(27, 13).Value = -238
$ws.Cells.Item(27, 14).Value = -801.6

$ws.Cells.Item(132, 8).Value = 4391.206
$ws.Cells.Item(132, 9).Value = 4213.433
$ws.Cells.Item(132, 10).Value = 5724.5
$ws.Cells.Item(132, 11).Value = 12640.299
$ws.Cells.Item(132, 12).Value = 17173.5
$ws.Cells.Item(132, 13).Value = -10110.299
$ws.Cells.Item(132, 14).Value = -22233.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 413.95456
$ws.Cells.Item(107, 9).Value = 193.42857
$ws.Cells.Item(107, 10).Value = 516.86664
$ws.Cells.Item(107, 11).Value = 580.28571
$ws.Cells.Item(107, 12).Value = 1550.59992
$ws.Cells.Item(107, 13).Value = 1339.71429
$ws.Cells.Item(107, 14).Value = -5390.59992

$ws.Cells.Item(132, 8).Value = 4310.4893
$ws.Cells.Item(132, 9).Value = 4359.405
$ws.Cells.Item(132, 10).Value = 3899.6
$ws.Cells.Item(132, 11).Value = 13078.215
$ws.Cells.Item(132, 12).Value = 11698.8
$ws.Cells.Item(132, 13).Value = -10548.215
$ws.Cells.Item(132, 14).Value = -16758.8
